$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8666075139614975
$ws.Range("C2").Value = 0.8353175135487731
$ws.Range("D2").Value = 0.6970035400410381
$ws.Range("E2").Value = 0.2833392684642462
$ws.Range("G2").Value = 0.002678115400934622
$ws.Range("I2").Value = 5.715146398380767
$ws.Range("J2").Value = 0.1466694704149489
$ws.Range("K2").Value = 1.649774853180617
$ws.Range("B3").Value = 0.8681249729491185
$ws.Range("C3").Value = 0.8057964673882054
$ws.Range("D3").Value = 0.6889004766573805
$ws.Range("E3").Value = 0.2792393688894919
$ws.Range("G3").Value = 0.002684530220707585
$ws.Range("I3").Value = 5.382151246756905
$ws.Range("J3").Value = 0.1438630810630315
$ws.Range("K3").Value = 1.627423254682128
$ws.Range("B4").Value = 0.8702139311409951
$ws.Range("C4").Value = 0.7880142068328269
$ws.Range("D4").Value = 0.6843063218402676
$ws.Range("E4").Value = 0.2768845755685234
$ws.Range("G4").Value = 0.002688664770115384
$ws.Range("I4").Value = 5.177040593450783
$ws.Range("J4").Value = 0.1422295955840411
$ws.Range("K4").Value = 1.615216906400576
$ws.Range("B5").Value = 0.871354729828937
$ws.Range("C5").Value = 0.7808534242635687
$ws.Range("D5").Value = 0.6825296678792938
$ws.Range("E5").Value = 0.2759656735488107
$ws.Range("G5").Value = 0.002690399074928368
$ws.Range("I5").Value = 5.093282016610061
$ws.Range("J5").Value = 0.1415863492707956
$ws.Range("K5").Value = 1.610622450493565
$ws.Range("B6").Value = 0.8715616033363176
$ws.Range("C6").Value = 0.7796695303866272
$ws.Range("D6").Value = 0.682240414055201
$ws.Range("E6").Value = 0.2758155432478873
$ws.Range("G6").Value = 0.002690690046953119
$ws.Range("I6").Value = 5.079363127211934
$ws.Range("J6").Value = 0.141480888599375
$ws.Range("K6").Value = 1.609882425042883
$ws.Range("B7").Value = 0.870228145952666
$ws.Range("C7").Value = 0.7879172883505419
$ws.Range("D7").Value = 0.6842819750262663
$ws.Range("E7").Value = 0.2768720183964177
$ws.Range("G7").Value = 0.002688687959255764
$ws.Range("I7").Value = 5.175911716005459
$ws.Range("J7").Value = 0.1422208299626746
$ws.Range("K7").Value = 1.615153408899005
$ws.Range("B8").Value = 0.8668896549933152
$ws.Range("C8").Value = 0.8250668001232953
$ws.Range("D8").Value = 0.6941302798307731
$ws.Range("E8").Value = 0.2818917704560278
$ws.Range("G8").Value = 0.002680286708672217
$ws.Range("I8").Value = 5.60045800844722
$ws.Range("J8").Value = 0.145683139980612
$ws.Range("K8").Value = 1.641751875347524
$ws.Range("B9").Value = 0.8695945567912702
$ws.Range("C9").Value = 0.9006863102993634
$ws.Range("D9").Value = 0.7164849389782546
$ws.Range("E9").Value = 0.293034620697064
$ws.Range("G9").Value = 0.002665356450182352
$ws.Range("I9").Value = 6.428359953805028
$ws.Range("J9").Value = 0.1531905914738587
$ws.Range("K9").Value = 1.706046155298111
$ws.Range("B10").Value = 0.8773234035767246
$ws.Range("C10").Value = 0.9579972347166859
$ws.Range("D10").Value = 0.7347915734137871
$ws.Range("E10").Value = 0.3020276605219081
$ws.Range("G10").Value = 0.002655315834802678
$ws.Range("I10").Value = 7.034584969946252
$ws.Range("J10").Value = 0.159153970347802
$ws.Range("K10").Value = 1.760821337617585
$ws.Range("B11").Value = 0.882109067168102
$ws.Range("C11").Value = 0.9844649502055631
$ws.Range("D11").Value = 0.7435350156981997
$ws.Range("E11").Value = 0.3062972261427745
$ws.Range("G11").Value = 0.002650946956166598
$ws.Range("I11").Value = 7.310106922087243
$ws.Range("J11").Value = 0.1619663494312533
$ws.Range("K11").Value = 1.787408796128517
$ws.Range("B12").Value = 0.8841057461811772
$ws.Range("C12").Value = 0.9945457328165048
$ws.Range("D12").Value = 0.7469061955026461
$ws.Range("E12").Value = 0.3079399331358061
$ws.Range("G12").Value = 0.002649320924890626
$ws.Range("I12").Value = 7.414417000860738
$ws.Range("J12").Value = 0.1630458236218857
$ws.Range("K12").Value = 1.79771944223512
$ws.Range("B13").Value = 0.8836674919096765
$ws.Range("C13").Value = 0.9923720588308811
$ws.Range("D13").Value = 0.7461774655297404
$ws.Range("E13").Value = 0.3075849903972596
$ws.Range("G13").Value = 0.00264966986107071
$ws.Range("I13").Value = 7.391952833564972
$ws.Range("J13").Value = 0.1628126925598963
$ws.Range("K13").Value = 1.795488033302576
$ws.Range("B14").Value = 0.882269627625675
$ws.Range("C14").Value = 0.9852931334204982
$ws.Range("D14").Value = 0.7438111549629411
$ws.Range("E14").Value = 0.3064318518873819
$ws.Range("G14").Value = 0.002650812613884499
$ws.Range("I14").Value = 7.318688987138046
$ws.Range("J14").Value = 0.1620548671213413
$ws.Range("K14").Value = 1.788252183988988
$ws.Range("B15").Value = 0.8814374725554615
$ws.Range("C15").Value = 0.9809646746768976
$ws.Range("D15").Value = 0.7423695786784208
$ws.Range("E15").Value = 0.3057289032594213
$ws.Range("G15").Value = 0.002651516271700635
$ws.Range("I15").Value = 7.27380999416755
$ws.Range("J15").Value = 0.1615925687668067
$ws.Range("K15").Value = 1.783851677746554
$ws.Range("B16").Value = 0.8770363581663219
$ws.Range("C16").Value = 0.9562755736632766
$ws.Range("D16").Value = 0.7342285717370771
$ws.Range("E16").Value = 0.3017522478784898
$ws.Range("G16").Value = 0.002655605331724922
$ws.Range("I16").Value = 7.016574916845713
$ws.Range("J16").Value = 0.1589721930246242
$ws.Range("K16").Value = 1.75911759327505
$ws.Range("B17").Value = 0.8746629441981497
$ws.Range("C17").Value = 0.9412319072945934
$ws.Range("D17").Value = 0.7293411085426271
$ws.Range("E17").Value = 0.2993586035928715
$ws.Range("G17").Value = 0.002658164576114101
$ws.Range("I17").Value = 6.858713416429197
$ws.Range("J17").Value = 0.157390304479847
$ws.Range("K17").Value = 1.744373445107584
$ws.Range("B18").Value = 0.8734172220201799
$ws.Range("C18").Value = 0.9326164656924334
$ws.Range("D18").Value = 0.7265690648868031
$ws.Range("E18").Value = 0.2979986376683286
$ws.Range("G18").Value = 0.002659655294903891
$ws.Range("I18").Value = 6.767891110906248
$ws.Range("J18").Value = 0.1564898095492708
$ws.Range("K18").Value = 1.736050036297684
$ws.Range("B19").Value = 0.8730158935381382
$ws.Range("C19").Value = 0.929705797903182
$ws.Range("D19").Value = 0.7256371996440976
$ws.Range("E19").Value = 0.297541054093557
$ws.Range("G19").Value = 0.002660163246319314
$ws.Range("I19").Value = 6.737135653359474
$ws.Range("J19").Value = 0.1561865210614144
$ws.Range("K19").Value = 1.733258771573645
$ws.Range("B20").Value = 0.8749032271207682
$ws.Range("C20").Value = 0.9428294661053087
$ws.Range("D20").Value = 0.7298573378112678
$ws.Range("E20").Value = 0.2996116712227632
$ws.Range("G20").Value = 0.002657890205087987
$ws.Range("I20").Value = 6.875520494976683
$ws.Range("J20").Value = 0.1575577287486425
$ws.Range("K20").Value = 1.74592671503035
$ws.Range("B21").Value = 0.8826751931983949
$ws.Range("C21").Value = 0.9873708023858399
$ws.Range("D21").Value = 0.74450455945464
$ws.Range("E21").Value = 0.3067698513848427
$ws.Range("G21").Value = 0.00265047619133112
$ws.Range("I21").Value = 7.340208910184856
$ws.Range("J21").Value = 0.1622770643149067
$ws.Range("K21").Value = 1.790370925923128
$ws.Range("B22").Value = 0.8888306493799121
$ws.Range("C22").Value = 1.016819880792127
$ws.Range("D22").Value = 0.7544286430703835
$ws.Range("E22").Value = 0.311599303374571
$ws.Range("G22").Value = 0.002645795965781777
$ws.Range("I22").Value = 7.643775001054848
$ws.Range("J22").Value = 0.1654459451010268
$ws.Range("K22").Value = 1.820832740681823
$ws.Range("B23").Value = 0.8854462822455957
$ws.Range("C23").Value = 1.001071024639486
$ws.Range("D23").Value = 0.7490996761642066
$ws.Range("E23").Value = 0.3090078227075566
$ws.Range("G23").Value = 0.002648278834074776
$ws.Range("I23").Value = 7.481764277409923
$ws.Range("J23").Value = 0.1637468647000304
$ws.Range("K23").Value = 1.804444411960986
$ws.Range("B24").Value = 0.8747942253624501
$ws.Range("C24").Value = 0.9421071059248334
$ws.Range("D24").Value = 0.7296238327793105
$ws.Range("E24").Value = 0.2994972089925767
$ws.Range("G24").Value = 0.002658014187804541
$ws.Range("I24").Value = 6.867922220099445
$ws.Range("J24").Value = 0.1574820083594233
$ws.Range("K24").Value = 1.745224004701129
$ws.Range("B25").Value = 0.8678620278080018
$ws.Range("C25").Value = 0.8799261168185808
$ws.Range("D25").Value = 0.7101088755126455
$ws.Range("E25").Value = 0.2898796387945239
$ws.Range("G25").Value = 0.002669231464149968
$ws.Range("I25").Value = 6.204805590287151
$ws.Range("J25").Value = 0.1510817531107804
$ws.Range("K25").Value = 1.687339449476866

Write-Output "Updated 192 cells (B,C,D,E,G,I,J,K for rows 2-25)"
